# Hybrid Strategy Injection + Scheduler Logic for Reply/Post
# Rewrites ACCOUNTS (sheet1), CONTENT_LINES (sheet2) and CALENDAR (sheet3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: ACCOUNTS
#   - drop tone_voice / bio / core_topics (old J:K:L)
#   - swap persona_type / content_lines column order (new H / I)
#   - real account identities injected into the first 5 rows
#   - placeholder rows renumbered/shifted down, new row 12 appended
# ---------------------------------------------------------------------
$wsAccounts = $wb.Worksheets.Item("ACCOUNTS")

# header: content_lines/persona_type swap to persona_type/content_lines
$wsAccounts.Cells.Item(1, 8).Value = "persona_type"
$wsAccounts.Cells.Item(1, 9).Value = "content_lines"

# only the cells that actually change per row (A/C/D identity + G status +
# H/I persona columns); auth_token/proxy (E/F) and platform (B) stay as-is
$wsAccounts.Cells.Item(2, 1).Value = "acc_samuel"
$wsAccounts.Cells.Item(2, 3).Value = "Samuel_MendozCD"
$wsAccounts.Cells.Item(2, 4).Value = "febrero202627"
$wsAccounts.Cells.Item(2, 8).Value = "policy_analyst"
$wsAccounts.Cells.Item(2, 9).Value = "politics, urbanism"

$wsAccounts.Cells.Item(3, 1).Value = "acc_mariate"
$wsAccounts.Cells.Item(3, 3).Value = "mariatemonto"
$wsAccounts.Cells.Item(3, 4).Value = "febrero202628"
$wsAccounts.Cells.Item(3, 8).Value = "coffee_snob"
$wsAccounts.Cells.Item(3, 9).Value = "lifestyle, culture"

$wsAccounts.Cells.Item(4, 1).Value = "acc_daniel"
$wsAccounts.Cells.Item(4, 3).Value = "Daniel_VargasCc"
$wsAccounts.Cells.Item(4, 4).Value = "Habiaunavez205@"
$wsAccounts.Cells.Item(4, 7).Value = "active"
$wsAccounts.Cells.Item(4, 8).Value = "tech_visionary"
$wsAccounts.Cells.Item(4, 9).Value = "tech, future"

$wsAccounts.Cells.Item(5, 1).Value = "acc_nguerrero"
$wsAccounts.Cells.Item(5, 3).Value = "NGuerrero16814"
$wsAccounts.Cells.Item(5, 4).Value = "Habiaunavez205@"
$wsAccounts.Cells.Item(5, 7).Value = "active"
$wsAccounts.Cells.Item(5, 8).Value = "shitposter"
$wsAccounts.Cells.Item(5, 9).Value = "memes, rant"

$wsAccounts.Cells.Item(6, 1).Value = "acc_revistavoces"
$wsAccounts.Cells.Item(6, 3).Value = "RevistavocesD"
$wsAccounts.Cells.Item(6, 4).Value = "Febrero202630"
$wsAccounts.Cells.Item(6, 7).Value = "active"
$wsAccounts.Cells.Item(6, 8).Value = "news_outlet"
$wsAccounts.Cells.Item(6, 9).Value = "news, headlines"

$wsAccounts.Cells.Item(7, 1).Value = "account_05"
$wsAccounts.Cells.Item(7, 3).Value = "user_placeholder_5"
$wsAccounts.Cells.Item(7, 9).Value = "general"

$wsAccounts.Cells.Item(8, 1).Value = "account_06"
$wsAccounts.Cells.Item(8, 3).Value = "user_placeholder_6"
$wsAccounts.Cells.Item(8, 9).Value = "general"

$wsAccounts.Cells.Item(9, 1).Value = "account_07"
$wsAccounts.Cells.Item(9, 3).Value = "user_placeholder_7"
$wsAccounts.Cells.Item(9, 9).Value = "general"

$wsAccounts.Cells.Item(10, 1).Value = "account_08"
$wsAccounts.Cells.Item(10, 3).Value = "user_placeholder_8"
$wsAccounts.Cells.Item(10, 9).Value = "general"

$wsAccounts.Cells.Item(11, 1).Value = "account_09"
$wsAccounts.Cells.Item(11, 3).Value = "user_placeholder_9"
$wsAccounts.Cells.Item(11, 9).Value = "general"

# brand-new row 12
$wsAccounts.Cells.Item(12, 1).Value = "account_10"
$wsAccounts.Cells.Item(12, 2).Value = "twitter"
$wsAccounts.Cells.Item(12, 3).Value = "user_placeholder_10"
$wsAccounts.Cells.Item(12, 4).Value = "password_here"
$wsAccounts.Cells.Item(12, 7).Value = "inactive"
$wsAccounts.Cells.Item(12, 8).Value = "general"
$wsAccounts.Cells.Item(12, 9).Value = "general"

# the old sheet had 12 columns (through L); remove the now-unused tail
# (J: tone_voice, K: bio, L: core_topics) so the sheet is 9 columns wide
$wsAccounts.Range("J1:L11").Clear()

# ---------------------------------------------------------------------
# Sheet 2: CONTENT_LINES
#   - same shape (A1:E4), new persona-aligned content lines
# ---------------------------------------------------------------------
$wsContent = $wb.Worksheets.Item("CONTENT_LINES")

$wsContent.Cells.Item(2, 2).Value = "General updates."
$wsContent.Cells.Item(2, 5).Value = "Public"

$wsContent.Cells.Item(3, 1).Value = "politics"
$wsContent.Cells.Item(3, 2).Value = "Political analysis."
$wsContent.Cells.Item(3, 3).Value = "Analyze the impact of {{topic}} on local policy."
$wsContent.Cells.Item(3, 4).Value = "daily"
$wsContent.Cells.Item(3, 5).Value = "Voters"

$wsContent.Cells.Item(4, 1).Value = "tech"
$wsContent.Cells.Item(4, 2).Value = "Tech trends."
$wsContent.Cells.Item(4, 3).Value = "Explain why {{topic}} matters."
$wsContent.Cells.Item(4, 4).Value = "mwf"
$wsContent.Cells.Item(4, 5).Value = "Techies"

# ---------------------------------------------------------------------
# Sheet 3: CALENDAR
#   - header: media_path -> target_url, hashtags -> action_type
#   - rows 2-4 replaced with "direct_*" reply posts, rows 5-11 appended
# ---------------------------------------------------------------------
$wsCalendar = $wb.Worksheets.Item("CALENDAR")

$wsCalendar.Cells.Item(1, 7).Value = "target_url"
$wsCalendar.Cells.Item(1, 8).Value = "action_type"

$calendarRows = @(
    @("direct_acc_samuel","acc_samuel","support_lgv","2026-02-14 11:50","approved","Es fundamental entender que la seguridad jurídica es la base de la inversión. Gran análisis @LuisGuillermoVl.","https://x.com/luisguillermovl/status/2022646985677840818","reply"),
    @("direct_acc_mariate","acc_mariate","support_lgv","2026-02-14 11:53","approved","Gracias por explicarlo tan claro. En las calles se siente la incertidumbre económica.","https://x.com/luisguillermovl/status/2022646985677840818","reply"),
    @("direct_acc_daniel","acc_daniel","support_lgv","2026-02-14 11:56","approved","¿Y dónde están los datos técnicos del gobierno? Vacíos. Gracias Concejal por poner los números sobre la mesa.","https://x.com/luisguillermovl/status/2022646985677840818","reply"),
    @("direct_acc_nguerrero","acc_nguerrero","support_lgv","2026-02-14 11:59","approved","Al fin alguien lo dice sin miedo. 🔥","https://x.com/luisguillermovl/status/2022646985677840818","reply"),
    @("direct_acc_revistavoces","acc_revistavoces","support_lgv","2026-02-14 12:02","approved","Claves del pronunciamiento del Concejal Vélez sobre el decreto de Salario Mínimo. Hilo 👇","https://x.com/luisguillermovl/status/2022646985677840818","quote"),
    @("debate_acc_samuel","acc_samuel","own_topic","2026-02-14 12:20","approved","La decisión del Consejo de Estado blinda nuestras instituciones. No es un capricho político.","","post"),
    @("debate_acc_mariate","acc_mariate","own_topic","2026-02-14 12:25","approved","Me preocupa mucho el costo de vida. ¿Qué piensan ustedes de este nuevo decreto?","","post"),
    @("debate_acc_daniel","acc_daniel","own_topic","2026-02-14 12:30","approved","Analizando el impacto en PYMES del nuevo decreto: Es insostenible sin subsidios cruzados.","","post"),
    @("debate_acc_nguerrero","acc_nguerrero","own_topic","2026-02-14 12:35","approved","El gobierno cree que somos tontos. Nos meten la mano al bolsillo y dicen que es 'justicia social'.","","post"),
    @("debate_acc_revistavoces","acc_revistavoces","own_topic","2026-02-14 12:40","approved","URGENTE: Reacciones encontradas tras la suspensión del decreto de salario mínimo. ¿Crisis institucional?","","post")
)

for ($r = 0; $r -lt $calendarRows.Count; $r++) {
    $rowData = $calendarRows[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $wsCalendar.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
